$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Just Visiting Jail" spacer row (A11:B11) carries a style that
# duplicates (save for an unused "apply fill" flag) the plain style already
# used elsewhere in the sheet (e.g. column N). Clearing its (already-none)
# fill pattern collapses it onto that shared style.
$ws.Range("A11:B11").Interior.Pattern = -4142

# Insert a new row 2 for the "START" board tile; everything below shifts
# down by one row.
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "START"
$ws.Range("B2").Value = "START"
$ws.Range("N2").Value = 0
$ws.Range("P2").Value = 7

# Move the selection onto the newly-added row, matching where the author
# left the cursor after adding it.
$ws.Range("P2").Select()
